# Applies the "add assign final asset weight for each security / remove sub
# category risk weight" edit across all three sheets (Alternative, Bond,
# Equity):
#   - Rename header C1 "Sub Category Risk Weight" -> "Sub Category Asset Weight"
#   - Rename header P1 "Asset Weight" -> "Portfolio Asset Weight"
#   - Update the Value at Risk 95% (N) figures with refreshed numbers
#   - Populate the (previously blank) Portfolio Asset Weight (P) column
#   - On the Equity sheet: drop the RWX row, fold the old VNQ/PDBC rows into
#     rows 7-8 (dropping the now-duplicate trailing rows), swap FLJP->FLJH and
#     FLKR->EWX for updated holdings, and refresh their figures.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Alternative sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Alternative")

$ws.Range("C1").Value = "Sub Category Asset Weight"
$ws.Range("P1").Value = "Portfolio Asset Weight"

$ws.Range("N2").Value = -10.42
$ws.Range("P2").Value = 7.11

$ws.Range("N3").Value = -4.47
$ws.Range("P3").Value = 12.89

$ws.Range("N4").Value = -9.18
$ws.Range("P4").Value = 0

$ws.Range("N5").Value = -80.28
$ws.Range("P5").Value = 0

$ws.Range("N6").Value = -8.72
$ws.Range("P6").Value = 0

$ws.Range("N7").Value = -7.56
$ws.Range("P7").Value = 0

# ---------------------------------------------------------------------------
# Bond sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Bond")

$ws.Range("C1").Value = "Sub Category Asset Weight"
$ws.Range("P1").Value = "Portfolio Asset Weight"

$ws.Range("N2").Value = -6.28
$ws.Range("P2").Value = 28.12

$ws.Range("N3").Value = 3.68
$ws.Range("P3").Value = 4.69

$ws.Range("N4").Value = -26.39
$ws.Range("P4").Value = 4.69

$ws.Range("N5").Value = -16.85
$ws.Range("P5").Value = 12.5

# ---------------------------------------------------------------------------
# Equity sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Equity")

# Drop the RWX / VNQ / PDBC rows (9-11): their content is superseded by the
# refreshed rows 7-8 below, and row 9 (RWX) is removed outright.
$ws.Rows("9:11").Delete()

$ws.Range("C1").Value = "Sub Category Asset Weight"
$ws.Range("P1").Value = "Portfolio Asset Weight"

# Row 2 - VOO
$ws.Range("C2").Value = 29.6
$ws.Range("N2").Value = -18.35
$ws.Range("P2").Value = 5.33

# Row 3 - FLCA
$ws.Range("C3").Value = 12.95
$ws.Range("N3").Value = -15.37
$ws.Range("P3").Value = 2.33

# Row 4 - FLJP -> FLJH (Franklin FTSE Japan Hedged ETF)
$ws.Range("A4").Value = "FLJH"
$ws.Range("C4").Value = 17.9
$ws.Range("D4").Value = "Franklin FTSE Japan Hedged ETF"
$ws.Range("I4").Value = 0.88
$ws.Range("J4").Value = 9.2
$ws.Range("K4").Value = 10.1
$ws.Range("L4").Value = 21.4
$ws.Range("M4").Value = 12.49
$ws.Range("N4").Value = -24.76
$ws.Range("O4").Value = 0.23
$ws.Range("P4").Value = 3.22

# Row 5 - FLAU
$ws.Range("C5").Value = 15.43
$ws.Range("N5").Value = -8.38
$ws.Range("P5").Value = 2.78

# Row 6 - FLKR -> EWX (SPDR S&P Emerging Markets Small Cap ETF)
$ws.Range("A6").Value = "EWX"
$ws.Range("C6").Value = 24.12
$ws.Range("D6").Value = "SPDR S&P Emerging Markets Small Cap ETF"
$ws.Range("E6").Value = "Diversified Emerging Mkts"
$ws.Range("H6").Value = 0.65
$ws.Range("I6").Value = 2.54
$ws.Range("J6").Value = 6.52
$ws.Range("K6").Value = 9.08
$ws.Range("L6").Value = 13.91
$ws.Range("M6").Value = 10.13
$ws.Range("N6").Value = -14.14
$ws.Range("O6").Value = 0.28
$ws.Range("P6").Value = 4.34

# Row 7 - SPEU -> VNQ (Vanguard Real Estate Index Fund)
$ws.Range("A7").Value = "VNQ"
$ws.Range("B7").Value = "REIT"
$ws.Range("C7").Value = 100
$ws.Range("D7").Value = "Vanguard Real Estate Index Fund"
$ws.Range("E7").Value = "Real Estate"
$ws.Range("H7").Value = 0.1
$ws.Range("I7").Value = 4.92
$ws.Range("J7").Value = 2.35
$ws.Range("K7").Value = 7.4
$ws.Range("L7").Value = 26.07
$ws.Range("M7").Value = 16.52
$ws.Range("N7").Value = -35.24
$ws.Range("O7").Value = 0.08
$ws.Range("P7").Value = 0

# Row 8 - SPEM -> PDBC (Invesco Optimum Yield Diversified Commodity Strategy No K-1 ETF)
$ws.Range("A8").Value = "PDBC"
$ws.Range("B8").Value = "Commodity"
$ws.Range("C8").Value = 100
$ws.Range("D8").Value = "Invesco Optimum Yield Diversified Commodity Strategy No K-1 ETF"
$ws.Range("E8").Value = "Commodities Broad Basket"
$ws.Range("F8").Value = "NasdaqGM"
$ws.Range("H8").Value = 0.59
$ws.Range("I8").Value = 12.91
$ws.Range("J8").Value = -1.11
$ws.Range("K8").Value = 11.83
$ws.Range("L8").Value = 7.96
$ws.Range("M8").Value = 9.12
$ws.Range("N8").Value = -0.58
$ws.Range("O8").Value = 0.83
$ws.Range("P8").Value = 12
